$wb = $excel.ActiveWorkbook

# Rename existing sheet
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "All_Results"

# Add a new worksheet right after "All_Results", for the best combination
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Best_Combination"

# Header row (same headers as All_Results)
$ws2.Range("A1").Value = "Chain"
$ws2.Range("B1").Value = "Start_Beta"
$ws2.Range("C1").Value = "Start_Alpha"
$ws2.Range("D1").Value = "Result_Alpha"
$ws2.Range("E1").Value = "Result_Beta"
$ws2.Range("F1").Value = "Pr_Chi"
$ws2.Range("G1").Value = "Deviance"

# Best combination row (copied from All_Results row 505: Start_Beta=170, Start_Alpha=200)
$ws2.Range("A2").Value = "Chain 1"
$ws2.Range("B2").Value = 170
$ws2.Range("C2").Value = 200

# D2/E2 hold numeric-looking text ("1911.672" / "121.5834"), same as the source
# data on All_Results - force them to stay text instead of being parsed as numbers.
$ws2.Range("D2:E2").NumberFormat = "@"
$ws2.Range("D2").Value = "1911.672"
$ws2.Range("E2").Value = "121.5834"
$ws2.Range("D2:E2").Style = "Normal"

$ws2.Range("F2").Value = 0.0000091614833760199
$ws2.Range("G2").Value = 23.2010049035514

$ws1.Select()
